$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

$ws.Range("A7").Value = "sure   bitti"
$ws.Range("A7").Select()
